$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 5-17 (Club, Participation Points, Performance Points,
# Total Points, Adjusted Total Points, ICL Eligible Number), re-sorted and
# updated per the fixed dictionary/matching for race validation.
$data = @(
    @(5,  "Newcastle Triathlon Club",      75, 319, 394, 300, 138),
    @(6,  "Balance Triathlon Club",        90, 238, 328, 300, 140),
    @(7,  "Brighton Baths Athletic Club",  90, 174, 264, 264, 36),
    @(8,  "Maitland Triathlon Club",       60, 182, 242, 242, 132),
    @(9,  "STG Triathlon Club",            90, 114, 204, 204, 49),
    @(11, "Tomaree Triathlon Club",        60, 65,  125, 125, 53),
    @(12, "Singleton Triathlon Club",      75, 46,  121, 121, 24),
    @(13, "Central Coast Triathlon Club",  45, 75,  120, 120, 135),
    @(14, "Concord Triathlon Club",        60, 50,  110, 110, 59),
    @(15, "Pulse Performance",             60, 36,  96,  96,  39),
    @(17, "Forster Triathlon Club",        60, 8,   68,  68,  46)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

$wb.Save()
